# "ajustando trataticas e referencias"
# Rename the single data sheet from the default "Sheet1" to "cns" so the
# tab/reference matches the rest of the workbook's naming (used elsewhere
# as e.g. cns!A1 / tabela_cns).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "cns"
